# Update project dates to 2025-08-01 to 2025-11-20 and format dates without time.
#
# Sheet "Sprint Backlog 3": column G (rows 4-93) holds "estimated completion"
# dates (3 rows per day). They get shifted forward from the Jan/Feb 2025
# range to start 2025-10-16, and change from the "yyyy-mm-dd h:mm:ss" style
# to a plain "yyyy-mm-dd" style (no time component).
#
# Sheet "SB BD Ch03": the burn-down chart date column A (and the paired
# Planned/Actual task counters in B/C) are recomputed for a longer,
# 36-day sprint window (2025-10-16 .. 2025-11-20), also switching to the
# "yyyy-mm-dd" style.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Sprint Backlog 3" -- column G estimated completion dates
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sprint Backlog 3")

$startSerial = 45946   # 2025-10-16, serial date number

for ($row = 4; $row -le 93; $row++) {
    $dayOffset = [Math]::Floor(($row - 4) / 3)
    $ws1.Cells.Item($row, 7).Value = $startSerial + $dayOffset
}

# Re-style the whole column at once (matches the single new cellXfs entry
# used for every G4:G93 cell in the target file).
$ws1.Range("G4:G93").NumberFormat = "yyyy-mm-dd"

# ---------------------------------------------------------------------
# Sheet 2: "SB BD Ch03" -- burn-down chart date / task counters
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SB BD Ch03")

$totalTasks = 90
$lastRowOffset = 35   # rows 3..38 -> 36 points, i.e. 35 steps down to 0

for ($row = 3; $row -le 38; $row++) {
    $i = $row - 3
    $ws2.Cells.Item($row, 1).Value = $startSerial + $i

    $remaining = [Math]::Round($totalTasks - ($totalTasks / $lastRowOffset) * $i, 1)
    $ws2.Cells.Item($row, 2).Value = $remaining
    $ws2.Cells.Item($row, 3).Value = $remaining
}

$ws2.Range("A3:A38").NumberFormat = "yyyy-mm-dd"
